# Handling of additional cases
# Insert a new mapping row ("series" -> Dose Unit -> dose_level_units) right
# after the existing "Dose" row, shifting all subsequent rows down by one,
# then refresh the ranges (dimension, autofilter, filter database) and the
# active selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92 (pushes the old rows 92-106 down to 93-107).
# The new row inherits formatting from the row above it (row 91), which
# already carries the highlighted style used for these "header mapping" rows.
$ws.Rows.Item(92).Insert()

$ws.Cells.Item(92, 1).Value = "series"
$ws.Cells.Item(92, 2).Value = "Dose Unit"
$ws.Cells.Item(92, 3).Value = "dose_level_units"

# Re-establish the autofilter over the new, larger range (A1:C107) without
# losing the filter header row.
$ws.AutoFilterMode = $False
[void]$ws.Range("A1:C107").AutoFilter()

# Update the (hidden) _FilterDatabase defined name so it also covers the
# newly expanded data range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$C`$107"
    }
}

# Move the active selection/view back to the top of the sheet.
[void]$ws.Range("B8").Select()
